$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: split the existing "Personal Leave" comment into an
#     Application cell (" Leave") and a Task cell ("Personal Leave"). ---
# Set D25 (Task) first so the "Personal Leave" shared string stays in
# use once C25 (Application) is overwritten with the new short text.
$ws.Range("D25").Value = "Personal Leave"
$ws.Range("C25").Value = " Leave"

# --- Row 32: fill in the previously-empty Application/Task cells. ---
$ws.Range("C32").Value = "Mujistore"
$ws.Range("D32").Value = "Fixing Mujistore issues"
# D32 already carried an explicit (bold-capable) cell style from the
# empty placeholder cell; touch the font so it normalizes to the same
# plain bordered style used by its neighboring data cells.
$ws.Range("D32").Font.Bold = $false

# --- Scroll / selection bookkeeping to match the saved view state. ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C25").Select()
